$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename school labels in column B (shared strings): drop the " ES" suffix,
# replace with two trailing spaces, for every row referencing that school.
$ws.Range("B2").Value = "Forest Grove  "
$ws.Range("B3").Value = "Forest Grove  "
$ws.Range("B4").Value = "Forest Grove  "
$ws.Range("B5").Value = "Forest Grove  "

$ws.Range("B6").Value = "Guilford  "
$ws.Range("B7").Value = "Guilford  "
$ws.Range("B8").Value = "Guilford  "
$ws.Range("B9").Value = "Guilford  "

$ws.Range("B10").Value = "Rolling Ridge  "
$ws.Range("B11").Value = "Rolling Ridge  "
$ws.Range("B12").Value = "Rolling Ridge  "
$ws.Range("B13").Value = "Rolling Ridge  "

$ws.Range("B14").Value = "Sterling  "
$ws.Range("B15").Value = "Sterling  "
$ws.Range("B16").Value = "Sterling  "
$ws.Range("B17").Value = "Sterling  "

$ws.Range("B18").Value = "Sugarland  "
$ws.Range("B19").Value = "Sugarland  "
$ws.Range("B20").Value = "Sugarland  "
$ws.Range("B21").Value = "Sugarland  "

$ws.Range("B22").Value = "Sully  "
$ws.Range("B23").Value = "Sully  "
$ws.Range("B24").Value = "Sully  "
$ws.Range("B25").Value = "Sully  "

# Update the Absence Rate (column D) values with the refreshed figures.
$ws.Range("D2").Value = 0.041
$ws.Range("D3").Value = 0.044
$ws.Range("D4").Value = 0.052
$ws.Range("D5").Value = 0.059
$ws.Range("D6").Value = 0.061
$ws.Range("D7").Value = 0.067
$ws.Range("D8").Value = 0.056
$ws.Range("D9").Value = 0.058
$ws.Range("D10").Value = 0.053
$ws.Range("D11").Value = 0.069
$ws.Range("D12").Value = 0.066
$ws.Range("D13").Value = 0.071
$ws.Range("D14").Value = 0.045
$ws.Range("D15").Value = 0.049
$ws.Range("D16").Value = 0.051
$ws.Range("D17").Value = 0.053
$ws.Range("D18").Value = 0.065
$ws.Range("D19").Value = 0.078
$ws.Range("D20").Value = 0.063
$ws.Range("D21").Value = 0.064
$ws.Range("D22").Value = 0.053
$ws.Range("D23").Value = 0.067
$ws.Range("D24").Value = 0.076
$ws.Range("D25").Value = 0.077

# Match the workbook's saved selection state (active cell D3).
$ws.Range("D3").Select()
